# Apply "custom accuracy" rounding to row 5 (B5:AH5) and remove row 6
# (the 1000-row sample was trimmed back down, dropping the last data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric measurements on row 5 to 2 decimal places (custom accuracy).
$row5 = @{
    "B5"  = 1.92
    "C5"  = 1.11
    "D5"  = 0.62
    "E5"  = 3.87
    "F5"  = 3.21
    "G5"  = 1.52
    "H5"  = 11.89
    "I5"  = 2.33
    "J5"  = 1
    "K5"  = 1.35
    "L5"  = 1.65
    "M5"  = 1.57
    "N5"  = 0.51
    "O5"  = 1.5
    "P5"  = 2.17
    "Q5"  = 1.48
    "R5"  = 0.71
    "S5"  = 0.28
    "T5"  = 15.67
    "U5"  = 4.63
    "V5"  = 1.39
    "W5"  = 3.01
    "X5"  = 1.66
    "Y5"  = 0.09
    "Z5"  = 5.17
    "AA5" = 1.23
    "AB5" = 1.23
    "AC5" = 1.41
    "AD5" = 1.65
    "AE5" = 0.5600000000000001
    "AF5" = 11.09
    "AG5" = 0.68
    "AH5" = 1.74
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Drop the last data row (row 6) entirely.
$ws.Rows.Item(6).Delete()
